# Extend the sheet with a new "2023" column (Q) mirroring the existing
# year columns (G:P), and fix an inconsistent capitalisation in column E
# ("Alternative Energy [ALT]" -> "Alternative energy [ALT]").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell Q1 ---------------------------------------------------
# Copy the look (bold / centered / bordered) of the neighbouring 2022
# header so the new header matches the rest of the row, then make it a
# genuine text value (like the other year headers) before typing "2023".
$ws.Range("P1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Q1").NumberFormat = "@"
$ws.Range("Q1").Value = "2023"

# --- Data rows for column Q (year 2023) -------------------------------
$ws.Range("Q2").Value = 3772
$ws.Range("Q3").Value = 30702
$ws.Range("Q4").Value = 0
$ws.Range("Q5").Value = 1939
$ws.Range("Q6").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("Q8").Value = 3719
$ws.Range("Q9").Value = 25480
$ws.Range("Q10").Value = 10
$ws.Range("Q11").Value = 2
$ws.Range("Q12").Value = 43
$ws.Range("Q13").Value = 3265
$ws.Range("Q14").Value = 53
$ws.Range("Q15").Value = 3283
$ws.Range("Q16").Value = 0
$ws.Range("Q17").Value = 1934
$ws.Range("Q18").Value = 0
$ws.Range("Q19").Value = 5
$ws.Range("Q20").Value = 3719
$ws.Range("Q21").Value = 25123
$ws.Range("Q22").Value = 0
$ws.Range("Q23").Value = 357
$ws.Range("Q24").Value = 0
$ws.Range("Q25").Value = 0
$ws.Range("Q26").Value = 0
$ws.Range("Q27").Value = 0

# Rows 28-31 (CNG / LNG beyond 2018) have no reported data in this
# dataset - keep Q blank there too, same as the existing blank cells in
# columns M:P for those rows. Copy the (unstyled) blank format from the
# neighbouring P cell so the cell exists but stays empty.
$ws.Range("P28").Copy()
$ws.Range("Q28").PasteSpecial(-4122)
$ws.Range("P29").Copy()
$ws.Range("Q29").PasteSpecial(-4122)
$ws.Range("P30").Copy()
$ws.Range("Q30").PasteSpecial(-4122)
$ws.Range("P31").Copy()
$ws.Range("Q31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("Q32").Value = 0
$ws.Range("Q33").Value = 16

# --- Fix inconsistent capitalisation of the "Alternative energy" label
$ws.Range("E14").Value = "Alternative energy [ALT]"
$ws.Range("E15").Value = "Alternative energy [ALT]"
